# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The "Periodo Mora" rows (16-22) are re-sorted into reverse chronological
# order (2112 .. 2106) and the "Salario Basico" column is refreshed for
# every worker/period line (908526 -> 738000). The "Valor Mora" figure of
# 19382 (previously tied to period 2112 in row 22) now belongs to the new
# first row (period 2112, row 16); every other period keeps 36341.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New period order (descending) for rows 16..22, with their refreshed
# "Valor Mora" (F) and "Salario Basico" (G) figures.
$periodos = @("2112", "2111", "2110", "2109", "2108", "2107", "2106")
$valorMora = @(19382, 36341, 36341, 36341, 36341, 36341, 36341)
$salarioBasico = @(738000, 738000, 738000, 738000, 738000, 738000, 738000)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
    $ws.Range("G$row").Value = $salarioBasico[$i]
}
